# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy style from an existing header cell so the
# new headers match the bold/bordered/centered look of the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows (2-45): every player on this roster shares the same team
# season record, so the same three numbers repeat down every row.
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 29).Value = 77   # AC - Wins
    $ws.Cells.Item($row, 30).Value = 85   # AD - Losses
    $ws.Cells.Item($row, 31).Value = 0    # AE - Ties
}
